# Regenerate save_data to use K (strikeouts landed, column G) instead of Strike#.
# Recalculated std/mean upstream produced new K values for each row; write them here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new K (column G) value
$newValues = @{
    2  = 2
    3  = 2
    4  = 3
    5  = 1
    6  = 3
    7  = 6
    8  = 6
    9  = 3
    10 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
